$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7955796666666667
$ws.Range("H2").Value = 2.386739
$ws.Range("I2").Value = 0.1186174580157865
$ws.Range("J2").Value = 0.1186174580157865
$ws.Range("M2").Value = 11.25749966666667
$ws.Range("N2").Value = 33.772499
$ws.Range("O2").Value = 0.6929800609896341
$ws.Range("P2").Value = 0.6929800609896341
$ws.Range("Q2").Value = 8.956237832306778
$ws.Range("R2").Value = 80.60614049076099
$ws.Range("S2").Value = 0.08219953329021512
$ws.Range("T2").Value = 0.08219953329021511
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7955796666666667
$ws.Range("H3").Value = 2.386739
$ws.Range("I3").Value = 0.1186174580157865
$ws.Range("J3").Value = 0.1186174580157865
$ws.Range("M3").Value = 0.9898276666666668
$ws.Range("O3").Value = 0.06093101107050686
$ws.Range("P3").Value = 0.06093101107050686
$ws.Range("Q3").Value = 0.7874867651041112
$ws.Range("R3").Value = 7.087380885937001
$ws.Range("S3").Value = 0.007227481647515272
$ws.Range("T3").Value = 0.007227481647515272
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7955796666666667
$ws.Range("H4").Value = 2.386739
$ws.Range("I4").Value = 0.1186174580157865
$ws.Range("J4").Value = 0.1186174580157865
$ws.Range("M4").Value = 3.821582
$ws.Range("N4").Value = 11.464746
$ws.Range("O4").Value = 0.2352458543950409
$ws.Range("P4").Value = 0.2352458543950409
$ws.Range("Q4").Value = 3.040372933699333
$ws.Range("R4").Value = 27.363356403294
$ws.Range("S4").Value = 0.0279042652570916
$ws.Range("T4").Value = 0.0279042652570916
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7955796666666667
$ws.Range("H5").Value = 2.386739
$ws.Range("I5").Value = 0.1186174580157865
$ws.Range("J5").Value = 0.1186174580157865
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1761463333333333
$ws.Range("N5").Value = 0.528439
$ws.Range("O5").Value = 0.01084307354481826
$ws.Range("P5").Value = 0.01084307354481827
$ws.Range("Q5").Value = 0.1401384411578889
$ws.Range("R5").Value = 1.261245970421
$ws.Range("S5").Value = 0.001286177820964566
$ws.Range("T5").Value = 0.001286177820964566
$ws.Range("I6").Value = 0.6312226244877757
$ws.Range("J6").Value = 0.6312226244877758
$ws.Range("M6").Value = 11.25749966666667
$ws.Range("N6").Value = 33.772499
$ws.Range("O6").Value = 0.6929800609896341
$ws.Range("P6").Value = 0.6929800609896341
$ws.Range("Q6").Value = 47.66060615877466
$ws.Range("R6").Value = 428.945455428972
$ws.Range("S6").Value = 0.4374246928155757
$ws.Range("T6").Value = 0.4374246928155757
$ws.Range("I7").Value = 0.6312226244877757
$ws.Range("J7").Value = 0.6312226244877758
$ws.Range("M7").Value = 0.9898276666666668
$ws.Range("O7").Value = 0.06093101107050686
$ws.Range("P7").Value = 0.06093101107050686
$ws.Range("R7").Value = 37.71548672852401
$ws.Range("S7").Value = 0.03846103272061905
$ws.Range("T7").Value = 0.03846103272061906
$ws.Range("I8").Value = 0.6312226244877757
$ws.Range("J8").Value = 0.6312226244877758
$ws.Range("M8").Value = 3.821582
$ws.Range("N8").Value = 11.464746
$ws.Range("O8").Value = 0.2352458543950409
$ws.Range("P8").Value = 0.2352458543950409
$ws.Range("Q8").Value = 16.179339995432
$ws.Range("R8").Value = 145.614059958888
$ws.Range("S8").Value = 0.1484925056111068
$ws.Range("T8").Value = 0.1484925056111069
$ws.Range("I9").Value = 0.6312226244877757
$ws.Range("J9").Value = 0.6312226244877758
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1761463333333333
$ws.Range("N9").Value = 0.528439
$ws.Range("O9").Value = 0.01084307354481826
$ws.Range("P9").Value = 0.01084307354481827
$ws.Range("Q9").Value = 0.7457465039213332
$ws.Range("R9").Value = 6.711718535292
$ws.Range("S9").Value = 0.006844393340474152
$ws.Range("T9").Value = 0.006844393340474156
$ws.Range("G10").Value = 1.666370333333333
$ws.Range("H10").Value = 4.999111
$ws.Range("I10").Value = 0.248448548064433
$ws.Range("J10").Value = 0.248448548064433
$ws.Range("M10").Value = 11.25749966666667
$ws.Range("N10").Value = 33.772499
$ws.Range("O10").Value = 0.6929800609896341
$ws.Range("P10").Value = 0.6929800609896341
$ws.Range("Q10").Value = 18.75916347204322
$ws.Range("R10").Value = 168.832471248389
$ws.Range("S10").Value = 0.1721698899904768
$ws.Range("T10").Value = 0.1721698899904768
$ws.Range("G11").Value = 1.666370333333333
$ws.Range("H11").Value = 4.999111
$ws.Range("I11").Value = 0.248448548064433
$ws.Range("J11").Value = 0.248448548064433
$ws.Range("M11").Value = 0.9898276666666668
$ws.Range("O11").Value = 0.06093101107050686
$ws.Range("P11").Value = 0.06093101107050686
$ws.Range("Q11").Value = 1.649419458845889
$ws.Range("R11").Value = 14.844775129613
$ws.Range("S11").Value = 0.01513822123256532
$ws.Range("T11").Value = 0.01513822123256532
$ws.Range("G12").Value = 1.666370333333333
$ws.Range("H12").Value = 4.999111
$ws.Range("I12").Value = 0.248448548064433
$ws.Range("J12").Value = 0.248448548064433
$ws.Range("M12").Value = 3.821582
$ws.Range("N12").Value = 11.464746
$ws.Range("O12").Value = 0.2352458543950409
$ws.Range("P12").Value = 0.2352458543950409
$ws.Range("Q12").Value = 6.368170871200666
$ws.Range("R12").Value = 57.313537840806
$ws.Range("S12").Value = 0.05844649096262491
$ws.Range("T12").Value = 0.05844649096262492
$ws.Range("G13").Value = 1.666370333333333
$ws.Range("H13").Value = 4.999111
$ws.Range("I13").Value = 0.248448548064433
$ws.Range("J13").Value = 0.248448548064433
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.1761463333333333
$ws.Range("N13").Value = 0.528439
$ws.Range("O13").Value = 0.01084307354481826
$ws.Range("P13").Value = 0.01084307354481827
$ws.Range("Q13").Value = 0.2935250241921111
$ws.Range("R13").Value = 2.641725217729
$ws.Range("S13").Value = 0.002693945878765961
$ws.Range("T13").Value = 0.002693945878765962
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.01147833333333333
$ws.Range("H14").Value = 0.034435
$ws.Range("I14").Value = 0.00171136943200476
$ws.Range("J14").Value = 0.00171136943200476
$ws.Range("M14").Value = 11.25749966666667
$ws.Range("N14").Value = 33.772499
$ws.Range("O14").Value = 0.6929800609896341
$ws.Range("P14").Value = 0.6929800609896341
$ws.Range("Q14").Value = 0.1292173336738889
$ws.Range("R14").Value = 1.162956003065
$ws.Range("S14").Value = 0.001185944893366454
$ws.Range("T14").Value = 0.001185944893366454
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.01147833333333333
$ws.Range("H15").Value = 0.034435
$ws.Range("I15").Value = 0.00171136943200476
$ws.Range("J15").Value = 0.00171136943200476
$ws.Range("M15").Value = 0.9898276666666668
$ws.Range("O15").Value = 0.06093101107050686
$ws.Range("P15").Value = 0.06093101107050686
$ws.Range("Q15").Value = 0.01136157190055556
$ws.Range("R15").Value = 0.102254147105
$ws.Range("S15").Value = 0.0001042754698072091
$ws.Range("T15").Value = 0.0001042754698072091
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.01147833333333333
$ws.Range("H16").Value = 0.034435
$ws.Range("I16").Value = 0.00171136943200476
$ws.Range("J16").Value = 0.00171136943200476
$ws.Range("M16").Value = 3.821582
$ws.Range("N16").Value = 11.464746
$ws.Range("O16").Value = 0.2352458543950409
$ws.Range("P16").Value = 0.2352458543950409
$ws.Range("Q16").Value = 0.04386539205666667
$ws.Range("R16").Value = 0.39478852851
$ws.Range("S16").Value = 0.0004025925642175156
$ws.Range("T16").Value = 0.0004025925642175157
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.01147833333333333
$ws.Range("H17").Value = 0.034435
$ws.Range("I17").Value = 0.00171136943200476
$ws.Range("J17").Value = 0.00171136943200476
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.1761463333333333
$ws.Range("N17").Value = 0.528439
$ws.Range("O17").Value = 0.01084307354481826
$ws.Range("P17").Value = 0.01084307354481827
$ws.Range("Q17").Value = 0.002021866329444444
$ws.Range("R17").Value = 0.018196796965
$ws.Range("S17").Value = 0.00001855650461358147
$ws.Range("T17").Value = 0.00001855650461358148
